$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57 (shifts existing rows 57:90 down to 58:91,
# and the used range grows from A1:R90 to A1:R91).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly price record.
$ws.Cells.Item(57, 1).Value = 1
$ws.Cells.Item(57, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(57, 4).Value = 44603
$ws.Cells.Item(57, 5).Value = 15
$ws.Cells.Item(57, 6).Value = 100112036
$ws.Cells.Item(57, 7).Value = "Caigua"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 130
$ws.Cells.Item(57, 11).Value = 10000
$ws.Cells.Item(57, 12).Value = 11000
$ws.Cells.Item(57, 13).Value = 10500
$ws.Cells.Item(57, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(57, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(57, 16).Value = 525
$ws.Cells.Item(57, 17).Value = 20
$ws.Cells.Item(57, 18).Value = "Hortaliza"
